# UPDATE technology portfolios for Norway
# - Update base investment cost (E2) on the "2025" sheet; the other year
#   sheets reference this value via formula and will recalc automatically.
# - Update COP (G2) on every year sheet to the new value.
# - Update charge/discharge cost (H2) on the "2025" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G2").Value = 8.5787499999999994
}

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("E2").Value = 2300
$ws2025.Range("H2").Value = 5

$excel.CalculateFullRebuild()
